$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 4 previously-blank rows (7-10) of the "Ark1" time log with the
# next batch of logged tasks for 2020-02-26. New values are entered in the
# same (slightly mixed) order Marc actually typed them in, which matters
# because it determines the order new entries land in the shared-string
# table: F7 first, then A7/B7, then column A top-to-bottom for rows 8-10,
# then B10, and finally column F for rows 10/8/9.

# --- Row 7 ---
$ws.Range("F7").Value = "30 min"
$ws.Range("A7").Value = "OC0802 angivBrugstidogScrapværdiogAnskaffelseværd"
$ws.Range("B7").Value = "System Analyst "
$ws.Range("C7").Value = "2/26/2020"
$ws.Range("D7").Value = 0.39583333333333331
$ws.Range("E7").Value = 0.4236111111111111

# --- Row 8 ---
$ws.Range("A8").Value = "Kundemøde KKO"
$ws.Range("C8").Value = "2/26/2020"
$ws.Range("D8").Value = 0.46875
$ws.Range("E8").Value = 0.49652777777777773

# --- Row 9 ---
$ws.Range("A9").Value = "Kundemøde Indtjeningsbidrag"
$ws.Range("C9").Value = "2/26/2020"
$ws.Range("D9").Value = 0.52777777777777779
$ws.Range("E9").Value = 0.55555555555555558

# --- Row 10 ---
$ws.Range("A10").Value = "OC06 beregn kontant kapacitetomkostning"
$ws.Range("B10").Value = "System Analyst "
$ws.Range("C10").Value = "2/26/2020"
$ws.Range("D10").Value = 0.5625
$ws.Range("E10").Value = 0.58680555555555558

$ws.Range("F10").Value = "30min"
$ws.Range("F8").Value = "45min"
$ws.Range("F9").Value = "45min"

# Column A widened (auto best-fit) after the longer strings were entered.
# (ColumnWidth set-value is rendered back out ~0.8333 wider than requested by
# this host's writer, so back the request off to land exactly on 59.)
$ws.Columns.Item(1).ColumnWidth = 58.166666666666664
